$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("network")
$ws.Rows(7).Insert()

$ws.Range("A7").Value = "MSN1"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 0
$ws.Range("V7").Value = 0
